$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Bitcoin)
$ws.Range("D2").Value = "52.325.21"

# Row 3 (Ethereum)
$ws.Range("D3").Value = "2.843.21"
$ws.Range("E3").Value = "  +1.41%  "

# Row 4 (TetherUSD)
$ws.Range("E4").Value = "  -0.04%  "

# Row 5 (BNB)
$ws.Range("E5").Value = "  +4.02%  "

# Row 6 (Solana)
$origStyle = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "112.65"
$ws.Range("D6").Style = $origStyle
$ws.Range("E6").Value = "  -2.71%  "

# Row 7 (XRP)
$origStyle = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.571"
$ws.Range("D7").Style = $origStyle

# Row 8 (USDC)
$ws.Range("E8").Value = "  -0.03%  "

# Row 9 (Cardano)
$origStyle = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.607"
$ws.Range("D9").Style = $origStyle
$ws.Range("E9").Value = "  +2.68%  "

# Row 10 (Avalanche)
$origStyle = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.14"
$ws.Range("D10").Style = $origStyle
$ws.Range("E10").Value = "  -2.87%  "

# Row 11 (Dogecoin)
$origStyle = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0871"
$ws.Range("D11").Style = $origStyle
$ws.Range("E11").Value = "  +1.00%  "

# Rows 12 & 13: Chainlink and TRON swap positions (TRON moves to row 12, Chainlink to row 13)
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$origStyle = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.132"
$ws.Range("D12").Style = $origStyle
$ws.Range("E12").Value = "  +0.92%  "

$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$origStyle = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.13"
$ws.Range("D13").Style = $origStyle
$ws.Range("E13").Value = "  +0.60%  "

# Row 14 (Polkadot)
$ws.Range("E14").Value = "  +0.06%  "

# Row 15 (WrappedliquidstakedEther2.0)
$ws.Range("D15").Value = "3.289.70"
$ws.Range("E15").Value = "  +1.23%  "

# Row 16 (WrappedEther)
$ws.Range("D16").Value = "2.819.70"
$ws.Range("E16").Value = "  +0.58%  "

# Row 17 (Polygon)
$origStyle = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.938"
$ws.Range("D17").Style = $origStyle
$ws.Range("E17").Value = "  +5.21%  "

# Row 18 (WrappedBTC)
$ws.Range("D18").Value = "52.226.95"
$ws.Range("E18").Value = "  +0.13%  "

# Row 19 (Uniswap)
$origStyle = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.59"
$ws.Range("D19").Style = $origStyle
$ws.Range("E19").Value = "  +4.17%  "

# Row 20 (ImmutableX)
$ws.Range("E20").Value = "  -0.83%  "

# Row 21 (InternetComputer(DFINITY))
$ws.Range("E21").Value = "  +0.98%  "

# Row 23 (BitcoinCash)
$origStyle = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "273.33"
$ws.Range("D23").Style = $origStyle
$ws.Range("E23").Value = "  +1.19%  "

# Row 24 (Litecoin)
$origStyle = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "70.61"
$ws.Range("D24").Style = $origStyle
$ws.Range("E24").Value = "  +0.69%  "

# Row 25 (PancakeSwap)
$ws.Range("E25").Value = "  +2.79%  "

# Row 26 (EthereumClassic)
$origStyle = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "27.13"
$ws.Range("D26").Style = $origStyle
$ws.Range("E26").Value = "  +1.03%  "

# Row 27 (Dai)
$ws.Range("E27").Value = "  +0.03%  "

# Row 28 (Cosmos)
$ws.Range("E28").Value = "  +1.09%  "

# Row 29 (Toncoin)
$ws.Range("E29").Value = "  -0.08%  "

# Row 30 (Kaspa)
$origStyle = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.143"
$ws.Range("D30").Style = $origStyle
$ws.Range("E30").Value = "  +1.98%  "

# Row 31 (VeChain)
$ws.Range("E31").Value = "  +8.56%  "

# Row 32 (InjectiveProtocol)
$origStyle = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "35.42"
$ws.Range("D32").Style = $origStyle
$ws.Range("E32").Value = "  +2.76%  "

# Row 33 (OKB)
$origStyle = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "52.51"
$ws.Range("D33").Style = $origStyle
$ws.Range("E33").Value = "  +4.32%  "

# Row 34 (Filecoin)
$ws.Range("E34").Value = "  +2.45%  "

# Row 35 (RenderToken)
$origStyle = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.58"
$ws.Range("D35").Style = $origStyle
$ws.Range("E35").Value = "  +12.93%  "

# Row 36 (Hedera)
$origStyle = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0854"
$ws.Range("D36").Style = $origStyle
$ws.Range("E36").Value = "  +2.35%  "

# Row 37 (FirstDigitalUSD)
$ws.Range("E37").Value = "  -0.07%  "

# Row 38 (LidoDAOToken)
$ws.Range("E38").Value = "  +2.40%  "

# Row 39 (ARBITRUM)
$ws.Range("E39").Value = "  -2.56%  "

# Row 40 (Celestia)
$origStyle = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.54"
$ws.Range("D40").Style = $origStyle
$ws.Range("E40").Value = "  -0.93%  "

# Row 41 (Stellar)
$ws.Range("E41").Value = "  +1.81%  "

# Row 42 (Stacks)
$ws.Range("E42").Value = "  -1.07%  "

# Row 43 (Monero)
$origStyle = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "126.66"
$ws.Range("D43").Style = $origStyle
$ws.Range("E43").Value = "  +0.18%  "

# Row 44 (EnergySwap)
$origStyle = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "23.18"
$ws.Range("D44").Style = $origStyle
$ws.Range("E44").Value = "  -0.76%  "

# Row 45 (WEMIXToken)
$ws.Range("E45").Value = "  -0.51%  "

# Row 46 (Maker)
$ws.Range("D46").Value = "2.095.75"
$ws.Range("E46").Value = "  +1.88%  "

# Row 47 (NEARProtocol)
$ws.Range("E47").Value = "  +1.34%  "

# Row 48 (ApeXProtocol)
$ws.Range("E48").Value = "  -1.22%  "

# Row 49 (THORChain)
$origStyle = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.90"
$ws.Range("D49").Style = $origStyle
$ws.Range("E49").Value = "  +5.63%  "

# Row 50 (SEI)
$origStyle = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.970"
$ws.Range("D50").Style = $origStyle
$ws.Range("E50").Value = "  +1.07%  "

# Row 51 (FraxShare)
$ws.Range("E51").Value = "  +2.98%  "

